$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 52, pushing the existing rows 52-55 down to 55-58
$ws.Range("A52:T54").EntireRow.Insert()

# Row 52 (new): August Red / Primera
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 44615
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100103
$ws.Range("H52").Value = "Frutos de hueso (carozo)"
$ws.Range("I52").Value = 100103006
$ws.Range("J52").Value = "Nectarín"
$ws.Range("K52").Value = "August Red"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 250
$ws.Range("N52").Value = 21000
$ws.Range("O52").Value = 22000
$ws.Range("P52").Value = 21500
$ws.Range("Q52").Value = "$/bandeja 18 kilos granel"
$ws.Range("R52").Value = "Región de O'Higgins"
$ws.Range("S52").Value = 1194
$ws.Range("T52").Value = 18

# Row 53 (new): June Pearl / Segunda
$ws.Range("A53").Value = 1
$ws.Range("B53").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C53").Value = "Arica y Parinacota"
$ws.Range("D53").Value = 44615
$ws.Range("E53").Value = 15
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100103
$ws.Range("H53").Value = "Frutos de hueso (carozo)"
$ws.Range("I53").Value = 100103006
$ws.Range("J53").Value = "Nectarín"
$ws.Range("K53").Value = "June Pearl"
$ws.Range("L53").Value = "Segunda"
$ws.Range("M53").Value = 300
$ws.Range("N53").Value = 21000
$ws.Range("O53").Value = 22000
$ws.Range("P53").Value = 21500
$ws.Range("Q53").Value = "$/bandeja 18 kilos granel"
$ws.Range("R53").Value = "Región de O'Higgins"
$ws.Range("S53").Value = 1194
$ws.Range("T53").Value = 18

# Row 54 (new): Venus / Segunda
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44615
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = "Frutos de hueso (carozo)"
$ws.Range("I54").Value = 100103006
$ws.Range("J54").Value = "Nectarín"
$ws.Range("K54").Value = "Venus"
$ws.Range("L54").Value = "Segunda"
$ws.Range("M54").Value = 300
$ws.Range("N54").Value = 20000
$ws.Range("O54").Value = 21000
$ws.Range("P54").Value = 20500
$ws.Range("Q54").Value = "$/bandeja 18 kilos granel"
$ws.Range("R54").Value = "Región de O'Higgins"
$ws.Range("S54").Value = 1139
$ws.Range("T54").Value = 18

$dim = $ws.UsedRange.Address()
"Final used range: $dim"
